function Set-Row($ws, $r, $a, $b, $c, $d, $e) {
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet LP1912: header updates ---
$ws1.Range("A2").Value = "Última actualización: 06:53:56"
$ws1.Range("A3").Value = "Total filas: 79"

# --- Sheet LP1912: data rows ---
Set-Row $ws1 8 "03:52:04" "04:46" "215A_EL PATO" 54 "LP1912"
Set-Row $ws1 10 "04:44:46" "04:46" "15_ABASTO" 2 "LP1912"
Set-Row $ws1 48 "06:53:56" "06:53" "17_ROMERO" 0 "LP1912"
Set-Row $ws1 49 "06:46:37" "06:54" "14_ABASTO" 8 "LP1912"
Set-Row $ws1 50 "06:53:56" "06:54" "17_ROMERO" 1 "LP1912"
Set-Row $ws1 51 "06:53:56" "07:03" "225_GOMEZ" 10 "LP1912"
Set-Row $ws1 52 "06:46:37" "07:04" "225_GOMEZ" 18 "LP1912"
Set-Row $ws1 53 "06:53:56" "07:06" "215C_EL PATO" 13 "LP1912"
Set-Row $ws1 54 "06:18:01" "07:07" "215C_EL PATO" 49 "LP1912"
Set-Row $ws1 55 "06:53:56" "07:13" "14X44_ABASTO" 20 "LP1912"
Set-Row $ws1 56 "06:18:01" "07:14" "14X44_ABASTO" 56 "LP1912"
Set-Row $ws1 57 "06:53:56" "07:20" "215A_EL PATO" 27 "LP1912"
Set-Row $ws1 58 "06:35:33" "07:21" "215A_EL PATO" 46 "LP1912"
Set-Row $ws1 59 "06:53:56" "07:23" "16_SANTA ANA" 30 "LP1912"
Set-Row $ws1 60 "06:46:37" "07:24" "16_SANTA ANA" 38 "LP1912"
Set-Row $ws1 61 "06:53:56" "07:28" "14_ABASTO" 35 "LP1912"
Set-Row $ws1 62 "06:46:37" "07:29" "14_ABASTO" 43 "LP1912"
Set-Row $ws1 63 "06:53:56" "07:33" "23_HERNANDEZ" 40 "LP1912"
Set-Row $ws1 64 "06:53:56" "07:35" "17X38_ROMERO" 42 "LP1912"
Set-Row $ws1 65 "06:53:56" "07:36" "27_EL RETIRO" 43 "LP1912"
Set-Row $ws1 66 "06:46:37" "07:36" "17X38_ROMERO" 50 "LP1912"
Set-Row $ws1 67 "06:46:37" "07:37" "27_EL RETIRO" 79 "LP1912"
Set-Row $ws1 68 "06:53:56" "07:43" "10_OLMOS" 50 "LP1912"
Set-Row $ws1 69 "06:18:01" "07:44" "10_OLMOS" 86 "LP1912"
Set-Row $ws1 70 "06:53:56" "07:49" "15_ABASTO" 56 "LP1912"
Set-Row $ws1 71 "06:35:33" "07:58" "23_HERNANDEZ" 83 "LP1912"
Set-Row $ws1 72 "06:53:56" "07:59" "23_HERNANDEZ" 66 "LP1912"
Set-Row $ws1 73 "06:53:56" "07:59" "11_ETCHEVERRY" 66 "LP1912"
Set-Row $ws1 74 "06:18:01" "08:00" "11_ETCHEVERRY" 102 "LP1912"
Set-Row $ws1 75 "06:46:37" "08:00" "23_HERNANDEZ" 74 "LP1912"
Set-Row $ws1 76 "06:53:56" "08:01" "16_SANTA ANA" 68 "LP1912"
Set-Row $ws1 77 "06:53:56" "08:03" "17X38_ROMERO" 70 "LP1912"
Set-Row $ws1 78 "06:53:56" "08:13" "10_OLMOS" 80 "LP1912"
Set-Row $ws1 79 "06:46:37" "08:14" "10_OLMOS" 88 "LP1912"
Set-Row $ws1 80 "06:53:56" "08:19" "17_ROMERO" 86 "LP1912"
Set-Row $ws1 81 "06:53:56" "08:33" "215C_EL PATO" 100 "LP1912"
Set-Row $ws1 82 "06:35:33" "08:34" "215C_EL PATO" 119 "LP1912"
Set-Row $ws1 83 "06:53:56" "08:47" "215A_EL PATO" 114 "LP1912"
Set-Row $ws1 84 "06:53:56" "08:51" "16_P MOR-SANTA ANA" 118 "LP1912"

# --- Sheet LP1912-215: header updates ---
$ws2.Range("A2").Value = "Última actualización: 06:53:56"
$ws2.Range("A3").Value = "Total filas: 15"

# --- Sheet LP1912-215: data rows ---
Set-Row $ws2 14 "06:53:56" "07:06" "215C_EL PATO" 13 "LP1912"
Set-Row $ws2 16 "06:53:56" "07:20" "215A_EL PATO" 27 "LP1912"
Set-Row $ws2 18 "06:53:56" "08:33" "215C_EL PATO" 100 "LP1912"
Set-Row $ws2 20 "06:53:56" "08:47" "215A_EL PATO" 114 "LP1912"

# --- Sheet 6203-6173: header updates ---
$ws3.Range("A2").Value = "Última actualización: 06:53:56"

# --- Sheet 6203-6173: data rows ---
Set-Row $ws3 7 "06:53:56" "07:27" "215A_LA PLATA" 34 "L6173"
Set-Row $ws3 8 "06:53:56" "08:09" "215A_LA PLATA" 76 "L6173"
Set-Row $ws3 10 "06:53:56" "08:22" "215C_LA PLATA" 89 "L6203"

Write-Output "Edit complete"
